$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths
$ws.Columns.Item(1).ColumnWidth = 23.6640625
$ws.Columns.Item(2).ColumnWidth = 24.33203125
$ws.Columns.Item(3).ColumnWidth = 29.6640625
$ws.Columns.Item(4).ColumnWidth = 21.6640625
$ws.Columns.Item(5).ColumnWidth = 18.33203125
$ws.Columns.Item(6).ColumnWidth = 18.33203125
$ws.Columns.Item(7).ColumnWidth = 38.83203125
$ws.Columns.Item(8).ColumnWidth = 40.83203125
$ws.Columns.Item(9).ColumnWidth = 70.33203125

# Header row
$ws.Cells.Item(1, 1).Value = "product_title"
$ws.Cells.Item(1, 2).Value = "product_description"
$ws.Cells.Item(1, 3).Value = "product_price"
$ws.Cells.Item(1, 4).Value = "product_category"
$ws.Cells.Item(1, 5).Value = "customizable"
$ws.Cells.Item(1, 6).Value = "occasion"
$ws.Cells.Item(1, 7).Value = "product_image"
$ws.Cells.Item(1, 8).Value = "text_mask"
$ws.Cells.Item(1, 9).Value = "modal_mask"

$imagePath = "/Users/leharbhatt/Desktop/samsung-galaxy-s10-black-back.png"

# Row 2
$ws.Cells.Item(2, 1).Value = "TestProduct1"
$ws.Cells.Item(2, 2).Value = "testing1"
$ws.Cells.Item(2, 3).Value = 29.79
$ws.Cells.Item(2, 4).Value = "Placard"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = "Christmas"
$ws.Cells.Item(2, 7).Value = $imagePath
$ws.Cells.Item(2, 8).Value = $imagePath
$ws.Cells.Item(2, 9).Value = $imagePath

# Row 3
$ws.Cells.Item(3, 1).Value = "TestProduct2"
$ws.Cells.Item(3, 2).Value = "testing2"
$ws.Cells.Item(3, 3).Value = 29.79
$ws.Cells.Item(3, 4).Value = "Placard"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = "Christmas"
$ws.Cells.Item(3, 7).Value = $imagePath
$ws.Cells.Item(3, 8).Value = $imagePath
$ws.Cells.Item(3, 9).Value = $imagePath

# Row 4
$ws.Cells.Item(4, 1).Value = "TestProduct3"
$ws.Cells.Item(4, 2).Value = "testing3"
$ws.Cells.Item(4, 3).Value = 29.79
$ws.Cells.Item(4, 4).Value = "Placard"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = "Christmas"
$ws.Cells.Item(4, 7).Value = $imagePath
$ws.Cells.Item(4, 8).Value = $imagePath
$ws.Cells.Item(4, 9).Value = $imagePath

# Apply Menlo font style to the image-path cells (G2:I4), one cell at a
# time so each cell's format resolves to the same single cellXfs slot
# (matches the 2-style workbook produced by the original edit).
for ($row = 2; $row -le 4; $row++) {
    for ($col = 7; $col -le 9; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Font.Name = "Menlo"
        $cell.Font.Size = 11
        $cell.Font.Color = 0
    }
}

# View settings
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("I1").Select()
